# Update "想去人数" (column F) figures across all four sheets to the
# freshly scraped values (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new F value
$updates = @{
    "展览" = @{
        2  = 2624
        5  = 1331
        9  = 183
        11 = 8264
        14 = 124
        16 = 266
        19 = 332
        20 = 10378
        25 = 384
        26 = 170
        30 = 2651
        31 = 2066
        36 = 4048
        38 = 2566
        39 = 2879
        40 = 1220
        41 = 151
        42 = 316
        43 = 252
        45 = 96
        46 = 94
        48 = 81
        49 = 59
    }
    "演出" = @{
        6  = 183
        7  = 42
        14 = 31
        18 = 25
    }
    "本地生活" = @{
        3 = 9
    }
    "全部类型" = @{
        2  = 2624
        3  = 183
        5  = 1331
        9  = 183
        11 = 8264
        14 = 124
        16 = 266
        19 = 332
        20 = 10379
        25 = 384
        26 = 170
        27 = 31
        31 = 2651
        32 = 2066
        35 = 4048
        37 = 2566
        38 = 2879
        40 = 1220
        41 = 151
        42 = 316
        43 = 252
        45 = 96
        46 = 94
        48 = 81
        49 = 59
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
